$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 1554.1
$ws.Range("I88").Value = 398.33334
$ws.Range("J88").Value = 2049.4285
$ws.Range("K88").Value = 398.33334
$ws.Range("L88").Value = 2049.4285
$ws.Range("M88").Value = 7.666659999999979
$ws.Range("N88").Value = -2861.4285
$ws.Range("H91").Value = 1554.1
$ws.Range("I91").Value = 398.33334
$ws.Range("J91").Value = 2049.4285
$ws.Range("K91").Value = 398.33334
$ws.Range("L91").Value = 2049.4285
$ws.Range("M91").Value = 1005.66666
$ws.Range("N91").Value = -4857.4285
$ws.Range("H101").Value = 900
$ws.Range("I101").Value = 900
$ws.Range("K101").Value = 2700
$ws.Range("M101").Value = -1078
$ws.Range("H111").Value = 1046
$ws.Range("I111").Value = 399
$ws.Range("J111").Value = 1477.3334
$ws.Range("K111").Value = 1197
$ws.Range("L111").Value = 4432.0002
$ws.Range("M111").Value = 1870
$ws.Range("N111").Value = -10566.0002
$ws.Range("H132").Value = 1728.5
$ws.Range("I132").Value = 1594.9412
$ws.Range("J132").Value = 3999
$ws.Range("K132").Value = 4784.8236
$ws.Range("L132").Value = 11997
$ws.Range("M132").Value = -2254.8236
$ws.Range("N132").Value = -17057
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("N133").ClearContents()
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 2563.5715
$ws.Range("I137").Value = 1435.6
$ws.Range("J137").Value = 3589
$ws.Range("K137").Value = 4306.799999999999
$ws.Range("L137").Value = 10767
$ws.Range("M137").Value = -1756.799999999999
$ws.Range("N137").Value = -15867
$ws.Range("H138").Value = 3758.9265
$ws.Range("I138").Value = 1125.3
$ws.Range("J138").Value = 4213
$ws.Range("K138").Value = 3375.9
$ws.Range("L138").Value = 12639
$ws.Range("M138").Value = 1764.1
$ws.Range("N138").Value = -22919

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 257012.25
$ws.Range("J6").Value = 9350
$ws.Range("L6").Value = 9350
$ws.Range("N6").Value = -9696
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("N51").ClearContents()
$ws.Range("H61").Value = 1402.0294
$ws.Range("I61").Value = 1402.0294
$ws.Range("K61").Value = 1402.0294
$ws.Range("M61").Value = -1190.0294
$ws.Range("H74").Value = 4380.684
$ws.Range("I74").Value = 1252.5454
$ws.Range("J74").Value = 8681.875
$ws.Range("K74").Value = 1252.5454
$ws.Range("L74").Value = 8681.875
$ws.Range("M74").Value = -378.5454
$ws.Range("N74").Value = -10429.875
$ws.Range("H77").Value = 4380.684
$ws.Range("I77").Value = 1252.5454
$ws.Range("J77").Value = 8681.875
$ws.Range("K77").Value = 6262.727
$ws.Range("L77").Value = 43409.375
$ws.Range("M77").Value = -1894.727
$ws.Range("N77").Value = -52145.375
$ws.Range("H88").Value = 543
$ws.Range("I88").Value = 1116.3334
$ws.Range("J88").Value = 256.33334
$ws.Range("K88").Value = 1116.3334
$ws.Range("L88").Value = 256.33334
$ws.Range("M88").Value = -710.3334
$ws.Range("N88").Value = -1068.33334
$ws.Range("H91").Value = 543
$ws.Range("I91").Value = 1116.3334
$ws.Range("J91").Value = 256.33334
$ws.Range("K91").Value = 1116.3334
$ws.Range("L91").Value = 256.33334
$ws.Range("M91").Value = 287.6666
$ws.Range("N91").Value = -3064.33334
$ws.Range("H97").Value = 546.1667
$ws.Range("I97").Value = 509.45456
$ws.Range("K97").Value = 509.45456
$ws.Range("M97").Value = -13.45456000000001
$ws.Range("H132").Value = 1096.0952
$ws.Range("I132").Value = 1100.9
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 3302.7
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -772.7000000000003
$ws.Range("N132").Value = -8060
$ws.Range("H136").Value = 1402.0294
$ws.Range("I136").Value = 1402.0294
$ws.Range("K136").Value = 4206.0882
$ws.Range("M136").Value = -1656.0882

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1604.125
$ws.Range("I20").Value = 1505.6666
$ws.Range("J20").Value = 1899.5
$ws.Range("K20").Value = 1505.6666
$ws.Range("L20").Value = 1899.5
$ws.Range("M20").Value = -1258.6666
$ws.Range("N20").Value = -2393.5
$ws.Range("H134").Value = 3245.4546
$ws.Range("I134").Value = 3161.9524
$ws.Range("K134").Value = 9485.8572
$ws.Range("M134").Value = -6950.8572

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4297.2
$ws.Range("I31").Value = 2154.7144
$ws.Range("J31").Value = 6171.875
$ws.Range("K31").Value = 2154.7144
$ws.Range("L31").Value = 6171.875
$ws.Range("M31").Value = -1859.7144
$ws.Range("N31").Value = -6761.875
$ws.Range("H34").Value = 4297.2
$ws.Range("I34").Value = 2154.7144
$ws.Range("J34").Value = 6171.875
$ws.Range("K34").Value = 2154.7144
$ws.Range("L34").Value = 6171.875
$ws.Range("M34").Value = -1952.7144
$ws.Range("N34").Value = -6575.875
$ws.Range("H62").Value = 70432.836
$ws.Range("J62").Value = 136666
$ws.Range("L62").Value = 136666
$ws.Range("N62").Value = -137914
$ws.Range("H65").Value = 70432.836
$ws.Range("J65").Value = 136666
$ws.Range("L65").Value = 683330
$ws.Range("N65").Value = -689570
$ws.Range("H132").Value = 2153.12
$ws.Range("I132").Value = 1905.7693
$ws.Range("J132").Value = 2421.0833
$ws.Range("K132").Value = 5717.3079
$ws.Range("L132").Value = 7263.249899999999
$ws.Range("M132").Value = -3187.3079
$ws.Range("N132").Value = -12323.2499
$ws.Range("H134").Value = 4274.6665
$ws.Range("I134").Value = 3456.6667
$ws.Range("J134").Value = 5092.6665
$ws.Range("K134").Value = 10370.0001
$ws.Range("L134").Value = 15277.9995
$ws.Range("M134").Value = -7835.000100000001
$ws.Range("N134").Value = -20347.9995
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 12187.5
$ws.Range("I56").Value = 12187.5
$ws.Range("K56").Value = 12187.5
$ws.Range("M56").Value = -11657.5
$ws.Range("H80").Value = 4513.6665
$ws.Range("I80").Value = 3194.3333
$ws.Range("J80").Value = 5833
$ws.Range("K80").Value = 9582.999899999999
$ws.Range("L80").Value = 17499
$ws.Range("M80").Value = -8646.999899999999
$ws.Range("N80").Value = -19371
$ws.Range("H83").Value = 4513.6665
$ws.Range("I83").Value = 3194.3333
$ws.Range("J83").Value = 5833
$ws.Range("K83").Value = 28748.9997
$ws.Range("L83").Value = 52497
$ws.Range("M83").Value = -24068.9997
$ws.Range("N83").Value = -61857

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3044.125
$ws.Range("J113").Value = 4999.75
$ws.Range("L113").Value = 4999.75
$ws.Range("N113").Value = -9339.75
$ws.Range("H132").Value = 2685.4546
$ws.Range("I132").Value = 2305.577
$ws.Range("J132").Value = 4096.4287
$ws.Range("K132").Value = 6916.731000000001
$ws.Range("L132").Value = 12289.2861
$ws.Range("M132").Value = -4386.731000000001
$ws.Range("N132").Value = -17349.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3916.5833
$ws.Range("I46").Value = 2833.1667
$ws.Range("K46").Value = 2833.1667
$ws.Range("M46").Value = -2645.1667
$ws.Range("H55").Value = 482.92
$ws.Range("I55").Value = 321.72223
$ws.Range("K55").Value = 321.72223
$ws.Range("M55").Value = -148.72223
$ws.Range("H61").Value = 2846.35
$ws.Range("I61").Value = 2495.3125
$ws.Range("J61").Value = 4250.5
$ws.Range("K61").Value = 2495.3125
$ws.Range("L61").Value = 4250.5
$ws.Range("M61").Value = -2293.3125
$ws.Range("N61").Value = -4654.5
$ws.Range("H82").Value = 1475.6296
$ws.Range("J82").Value = 1177.2858
$ws.Range("L82").Value = 1177.2858
$ws.Range("N82").Value = -1899.2858
$ws.Range("H85").Value = 1475.6296
$ws.Range("J85").Value = 1177.2858
$ws.Range("L85").Value = 1177.2858
$ws.Range("N85").Value = -3673.2858
$ws.Range("H113").Value = 2846.35
$ws.Range("I113").Value = 2495.3125
$ws.Range("J113").Value = 4250.5
$ws.Range("K113").Value = 2495.3125
$ws.Range("L113").Value = 4250.5
$ws.Range("M113").Value = -325.3125
$ws.Range("N113").Value = -8590.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1142
$ws.Range("I81").Value = 1260.8182
$ws.Range("J81").Value = 955.2857
$ws.Range("K81").Value = 2521.6364
$ws.Range("L81").Value = 1910.5714
$ws.Range("M81").Value = -1460.6364
$ws.Range("N81").Value = -4032.5714
$ws.Range("H84").Value = 1142
$ws.Range("I84").Value = 1260.8182
$ws.Range("J84").Value = 955.2857
$ws.Range("K84").Value = 12608.182
$ws.Range("L84").Value = 9552.857
$ws.Range("M84").Value = -7304.181999999999
$ws.Range("N84").Value = -20160.857
$ws.Range("H132").Value = 1338.3334
$ws.Range("I132").Value = 1353.2
$ws.Range("K132").Value = 4059.6
$ws.Range("M132").Value = -1529.6
